# Daily attendance processing - 2025-11-29 01:24:14
#
# Normalizes the order of names in the "Recorded By" column (G) for the
# specific rows touched by today's processing run: the first two
# comma-separated entries are swapped (any trailing "System" entry stays
# last).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToFix = @(
    2, 3, 6, 7, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26,
    28, 29, 32, 33, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50,
    52, 54, 55, 58, 59, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74,
    76, 78, 83, 84, 85, 86, 87, 90, 92, 93, 94, 96, 99, 101, 109, 110,
    111, 112, 113, 116, 118, 119, 120, 122, 125, 127, 135, 136, 137,
    138, 139, 142, 144, 145, 146, 148, 151, 153
)

foreach ($r in $rowsToFix) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    $parts = $current -split ', '

    if ($parts.Count -eq 2) {
        $updated = $parts[1] + ", " + $parts[0]
    } elseif ($parts.Count -eq 3) {
        $updated = $parts[1] + ", " + $parts[0] + ", " + $parts[2]
    } else {
        $updated = $current
    }

    $cell.Value = $updated
}
